# Scheduled runner update: refresh market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-leve profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2701
$ws.Range("I29").Value = 1551.5
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 4654.5
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -4373.5
$ws.Range("N29").Value = -15562
$ws.Range("H38").Value = 219.66667
$ws.Range("I38").Value = 103.27273
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 309.81819
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = 62.18181000000004
$ws.Range("N38").Value = -5244
$ws.Range("H58").Value = 3844.0293
$ws.Range("I58").Value = 120
$ws.Range("J58").Value = 4340.567
$ws.Range("K58").Value = 360
$ws.Range("L58").Value = 13021.701
$ws.Range("M58").Value = -210
$ws.Range("N58").Value = -13321.701
$ws.Range("H62").Value = 3251.25
$ws.Range("I62").Value = 3251.25
$ws.Range("K62").Value = 3251.25
$ws.Range("M62").Value = -2627.25
$ws.Range("H65").Value = 3251.25
$ws.Range("I65").Value = 3251.25
$ws.Range("K65").Value = 16256.25
$ws.Range("M65").Value = -13136.25
$ws.Range("H87").Value = 35225
$ws.Range("J87").Value = 35225
$ws.Range("L87").Value = 35225
$ws.Range("N87").Value = -37721
$ws.Range("H90").Value = 35225
$ws.Range("J90").Value = 35225
$ws.Range("L90").Value = 105675
$ws.Range("N90").Value = -118155
$ws.Range("H112").Value = 29241124
$ws.Range("J112").Value = 3269442
$ws.Range("L112").Value = 9808326
$ws.Range("N112").Value = -9810542
$ws.Range("H113").Value = 3745
$ws.Range("I113").Value = 3967.5
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 3967.5
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = -713.5
$ws.Range("N113").Value = -9808
$ws.Range("H129").Value = 886.175
$ws.Range("I129").Value = 540.3
$ws.Range("J129").Value = 1001.4667
$ws.Range("K129").Value = 1620.9
$ws.Range("L129").Value = 3004.4001
$ws.Range("M129").Value = 3379.1
$ws.Range("N129").Value = -13004.4001
$ws.Range("H138").Value = 2735029.2
$ws.Range("I138").Value = 1809.8948
$ws.Range("J138").Value = 3971485.8
$ws.Range("K138").Value = 5429.6844
$ws.Range("L138").Value = 11914457.4
$ws.Range("M138").Value = -289.6844000000001
$ws.Range("N138").Value = -11924737.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6928.8335
$ws.Range("I88").Value = 3840.5715
$ws.Range("J88").Value = 8894.091
$ws.Range("K88").Value = 3840.5715
$ws.Range("L88").Value = 8894.091
$ws.Range("M88").Value = -3434.5715
$ws.Range("N88").Value = -9706.091
$ws.Range("H91").Value = 6928.8335
$ws.Range("I91").Value = 3840.5715
$ws.Range("J91").Value = 8894.091
$ws.Range("K91").Value = 3840.5715
$ws.Range("L91").Value = 8894.091
$ws.Range("M91").Value = -2436.5715
$ws.Range("N91").Value = -11702.091

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10700.897
$ws.Range("I86").Value = 10430.148
$ws.Range("J86").Value = 11310.083
$ws.Range("K86").Value = 10430.148
$ws.Range("L86").Value = 11310.083
$ws.Range("M86").Value = -9307.147999999999
$ws.Range("N86").Value = -13556.083
$ws.Range("H89").Value = 10700.897
$ws.Range("I89").Value = 10430.148
$ws.Range("J89").Value = 11310.083
$ws.Range("K89").Value = 52150.74
$ws.Range("L89").Value = 56550.415
$ws.Range("M89").Value = -46534.74
$ws.Range("N89").Value = -67782.41500000001
$ws.Range("H141").Value = 41421.6
$ws.Range("J141").Value = 41896.668
$ws.Range("L141").Value = 41896.668
$ws.Range("N141").Value = -52256.668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 37038184
$ws.Range("I16").Value = 703.3333
$ws.Range("K16").Value = 703.3333
$ws.Range("M16").Value = -416.3333
$ws.Range("H31").Value = 1705.1351
$ws.Range("I31").Value = 993.6667
$ws.Range("J31").Value = 3626.1
$ws.Range("K31").Value = 993.6667
$ws.Range("L31").Value = 3626.1
$ws.Range("M31").Value = -698.6667
$ws.Range("N31").Value = -4216.1
$ws.Range("H34").Value = 1705.1351
$ws.Range("I34").Value = 993.6667
$ws.Range("J34").Value = 3626.1
$ws.Range("K34").Value = 993.6667
$ws.Range("L34").Value = 3626.1
$ws.Range("M34").Value = -791.6667
$ws.Range("N34").Value = -4030.1
$ws.Range("H52").Value = 48500
$ws.Range("J52").Value = 48500
$ws.Range("L52").Value = 48500
$ws.Range("N52").Value = -49088
$ws.Range("H113").Value = 37038184
$ws.Range("I113").Value = 703.3333
$ws.Range("K113").Value = 703.3333
$ws.Range("M113").Value = 1466.6667
$ws.Range("H132").Value = 24316.844
$ws.Range("I132").Value = 1716.6316
$ws.Range("J132").Value = 147003.72
$ws.Range("K132").Value = 5149.8948
$ws.Range("L132").Value = 441011.16
$ws.Range("M132").Value = -2619.8948
$ws.Range("N132").Value = -446071.16

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 47619788
$ws.Range("I12").Value = 142858750
$ws.Range("J12").Value = 303.7143
$ws.Range("K12").Value = 428576250
$ws.Range("L12").Value = 911.1428999999999
$ws.Range("M12").Value = -428576077
$ws.Range("N12").Value = -1257.1429
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 449
$ws.Range("I22").Value = 449
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1347
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1178
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 449
$ws.Range("I27").Value = 449
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1347
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1245
$ws.Range("N27").ClearContents()
$ws.Range("H131").Value = 968.931
$ws.Range("J131").Value = 1031.6884
$ws.Range("L131").Value = 3095.0652
$ws.Range("N131").Value = -13175.0652
$ws.Range("H136").Value = 1492.4445
$ws.Range("I136").Value = 1112.375
$ws.Range("K136").Value = 3337.125
$ws.Range("M136").Value = 1762.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 31754.236
$ws.Range("I70").Value = 50725.137
$ws.Range("K70").Value = 50725.137
$ws.Range("M70").Value = -50455.137
$ws.Range("H73").Value = 31754.236
$ws.Range("I73").Value = 50725.137
$ws.Range("K73").Value = 50725.137
$ws.Range("M73").Value = -49789.137
$ws.Range("H102").Value = 1040.1538
$ws.Range("I102").Value = 913.55554
$ws.Range("K102").Value = 913.55554
$ws.Range("M102").Value = 708.44446
$ws.Range("H135").Value = 31540
$ws.Range("J135").Value = 31540
$ws.Range("L135").Value = 31540
$ws.Range("N135").Value = -41680

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 49639.812
$ws.Range("I132").Value = 28851.158
$ws.Range("K132").Value = 86553.474
$ws.Range("M132").Value = -84023.474

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 40948
$ws.Range("I136").Value = 23563.863
$ws.Range("J136").Value = 168431.67
$ws.Range("K136").Value = 70691.58900000001
$ws.Range("L136").Value = 505295.01
$ws.Range("M136").Value = -68141.58900000001
$ws.Range("N136").Value = -510395.01
$ws.Range("H137").Value = 49800
$ws.Range("J137").Value = 49800
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000
$ws.Range("H140").Value = 52772.5
$ws.Range("J140").Value = 52772.5
$ws.Range("L140").Value = 52772.5
$ws.Range("N140").Value = -63132.5
$ws.Range("H141").Value = 48500
$ws.Range("J141").Value = 48500
$ws.Range("L141").Value = 48500
$ws.Range("N141").Value = -58860
